$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 13551
$ws.Range("E2").Value = 372
$ws.Range("F2").Value = 372
$ws.Range("G2").Value = 411
$ws.Range("H2").Value = 305
$ws.Range("I2").Value = 305
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 10814
$ws.Range("L2").Value = 5060
$ws.Range("M2").Value = 5754
$ws.Range("N2").Value = 5754
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 315
$ws.Range("Q2").Value = 604
$ws.Range("R2").Value = -482
$ws.Range("S2").Value = 421
$ws.Range("T2").Value = 302
$ws.Range("U2").Value = 302
$ws.Range("V2").Value = 2625
$ws.Range("W2").Value = 2.75
$ws.Range("X2").Value = 2.25
$ws.Range("Y2").Value = 6.06
$ws.Range("Z2").Value = 3.27
$ws.Range("AA2").Value = 87.93000000000001
$ws.Range("AB2").Value = 1711.21
$ws.Range("AC2").Value = 5779
$ws.Range("AD2").Value = 10.68
$ws.Range("AE2").Value = 96498
$ws.Range("AF2").Value = 0.64
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 2.43
$ws.Range("AI2").Value = 29.37
$ws.Range("AJ2").Value = 6136619

# Row 3
$ws.Range("D3").Value = 14389
$ws.Range("E3").Value = 758
$ws.Range("F3").Value = 758
$ws.Range("G3").Value = 665
$ws.Range("H3").Value = 510
$ws.Range("I3").Value = 475
$ws.Range("J3").Value = 35
$ws.Range("K3").Value = 14378
$ws.Range("L3").Value = 7345
$ws.Range("M3").Value = 7034
$ws.Range("N3").Value = 6078
$ws.Range("O3").Value = 956
$ws.Range("P3").Value = 315
$ws.Range("Q3").Value = 322
$ws.Range("R3").Value = -525
$ws.Range("S3").Value = -311
$ws.Range("T3").Value = 626
$ws.Range("U3").Value = -304
$ws.Range("V3").Value = 4634
$ws.Range("W3").Value = 5.27
$ws.Range("X3").Value = 3.55
$ws.Range("Y3").Value = 8.039999999999999
$ws.Range("Z3").Value = 4.05
$ws.Range("AA3").Value = 104.42
$ws.Range("AB3").Value = 1821.89
$ws.Range("AC3").Value = 7544
$ws.Range("AD3").Value = 13.72
$ws.Range("AE3").Value = 101930
$ws.Range("AF3").Value = 1.02
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 1.45
$ws.Range("AI3").Value = 18.83
$ws.Range("AJ3").Value = 6136619

# Row 4
$ws.Range("D4").Value = 20009
$ws.Range("E4").Value = 1473
$ws.Range("F4").Value = 1473
$ws.Range("G4").Value = 1362
$ws.Range("H4").Value = 1012
$ws.Range("I4").Value = 898
$ws.Range("J4").Value = 113
$ws.Range("K4").Value = 20638
$ws.Range("L4").Value = 8977
$ws.Range("M4").Value = 11660
$ws.Range("N4").Value = 10592
$ws.Range("O4").Value = 1068
$ws.Range("P4").Value = 533
$ws.Range("Q4").Value = 1796
$ws.Range("R4").Value = -93
$ws.Range("S4").Value = -441
$ws.Range("T4").Value = 862
$ws.Range("U4").Value = 934
$ws.Range("V4").Value = 5314
$ws.Range("W4").Value = 7.36
$ws.Range("X4").Value = 5.05
$ws.Range("Y4").Value = 10.78
$ws.Range("Z4").Value = 5.78
$ws.Range("AA4").Value = 76.98999999999999
$ws.Range("AB4").Value = 1944.63
$ws.Range("AC4").Value = 8456
$ws.Range("AD4").Value = 12.42
$ws.Range("AE4").Value = 105764
$ws.Range("AF4").Value = 0.99
$ws.Range("AG4").Value = 1500
$ws.Range("AH4").Value = 1.43
$ws.Range("AI4").Value = 16.74
$ws.Range("AJ4").Value = 10289803

# Row 5
$ws.Range("D5").Value = 20412
$ws.Range("E5").Value = 889
$ws.Range("F5").Value = 889
$ws.Range("G5").Value = 700
$ws.Range("H5").Value = 492
$ws.Range("I5").Value = 355
$ws.Range("J5").Value = 137
$ws.Range("K5").Value = 21431
$ws.Range("L5").Value = 8968
$ws.Range("M5").Value = 12463
$ws.Range("N5").Value = 10825
$ws.Range("O5").Value = 1540
$ws.Range("P5").Value = 533
$ws.Range("Q5").Value = -501
$ws.Range("R5").Value = -2123
$ws.Range("S5").Value = 1324
$ws.Range("T5").Value = 735
$ws.Range("U5").Value = -1236
$ws.Range("V5").Value = 5654
$ws.Range("W5").Value = 4.35
$ws.Range("X5").Value = 2.41
$ws.Range("Y5").Value = 3.31
$ws.Range("Z5").Value = 2.34
$ws.Range("AA5").Value = 71.95999999999999
$ws.Range("AB5").Value = 1982.5
$ws.Range("AC5").Value = 3330
$ws.Range("AD5").Value = 27.87
$ws.Range("AE5").Value = 108090
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 1250
$ws.Range("AH5").Value = 1.35
$ws.Range("AI5").Value = 35.32
$ws.Range("AJ5").Value = 10289803

# Row 6
$ws.Range("D6").Value = 21238
$ws.Range("E6").Value = 970
$ws.Range("F6").Value = 970
$ws.Range("G6").Value = 806
$ws.Range("H6").Value = 597
$ws.Range("I6").Value = 450
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = 21739
$ws.Range("L6").Value = 9079
$ws.Range("M6").Value = 12660
$ws.Range("N6").Value = 10947
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 533
$ws.Range("Q6").Value = 684
$ws.Range("R6").Value = -526
$ws.Range("S6").Value = -176
$ws.Range("T6").Value = 1403
$ws.Range("U6").Value = -719
$ws.Range("V6").Value = 5308
$ws.Range("W6").Value = 4.57
$ws.Range("X6").Value = 2.81
$ws.Range("Y6").Value = 4.13
$ws.Range("Z6").Value = 2.77
$ws.Range("AA6").Value = 71.72
$ws.Range("AB6").Value = 2032.28
$ws.Range("AC6").Value = 4219
$ws.Range("AD6").Value = 13.04
$ws.Range("AE6").Value = 109300
$ws.Range("AF6").Value = 0.5
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 33.44
$ws.Range("AJ6").Value = 10289803

# Row 7: clear all data columns, keep A/B/C
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all data columns, keep A/B/C
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all data columns, keep A/B/C
$ws.Range("D9:AJ9").ClearContents()
